$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6851942329906543
$ws.Range("D2").Value = 0.5003778768378944

$ws.Range("C3").Value = 1.444797119256442
$ws.Range("D3").Value = 0.1626082845614432

$ws.Range("C4").Value = 1.743474057542678
$ws.Range("D4").Value = 0.095211155056798

$ws.Range("C5").Value = 3.008550880283925
$ws.Range("D5").Value = 0.006465173062416785

$ws.Range("C6").Value = 0.7013207802777968
$ws.Range("D6").Value = 0.4904555296657263

$ws.Range("C7").Value = 1.234627775944518
$ws.Range("D7").Value = 0.2299941153296448

$ws.Range("C8").Value = 2.176680856983311
$ws.Range("D8").Value = 0.04051711651350454

$ws.Range("C9").Value = 0.09899821313638568
$ws.Range("D9").Value = 0.9220362152518293

$ws.Range("C10").Value = 1.114648576790221
$ws.Range("D10").Value = 0.2770302110246514

$ws.Range("C11").Value = 1.397624387595759
$ws.Range("D11").Value = 0.1761615727352215
